$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Force the value to be stored as text (avoids Excel auto-converting
    # numeric-looking strings like "577.86" into real numbers), while
    # clearing the explicit number-format style afterwards so the cell
    # keeps its original (default/general) style.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "62.605.37"
$ws.Range("E2").Value = "  -1.08%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.436.34"
$ws.Range("E3").Value = "  -1.52%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
Set-TextValue "D5" "577.86"
$ws.Range("E5").Value = "  -1.28%  "

# Row 6 - Solana
Set-TextValue "D6" "147.08"
$ws.Range("E6").Value = "  -0.77%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.05%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.29%  "

# Row 9 - Toncoin
$ws.Range("E9").Value = "  +3.47%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -2.93%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +2.02%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "4.025.81"
$ws.Range("E12").Value = "  -1.45%  "

# Row 14 - Avalanche
Set-TextValue "D14" "28.11"
$ws.Range("E14").Value = "  -5.85%  "

# Row 15 - WrappedEther
Set-TextValue "D15" "3.431.00"
$ws.Range("E15").Value = "  -1.67%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -1.71%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "62.660.63"
$ws.Range("E17").Value = "  -1.04%  "

# Row 18 - Polkadot
Set-TextValue "D18" "6.37"
$ws.Range("E18").Value = "  +0.49%  "

# Row 19 - Chainlink
Set-TextValue "D19" "14.55"
$ws.Range("E19").Value = "  +1.02%  "

# Row 20 - Uniswap
Set-TextValue "D20" "9.04"
$ws.Range("E20").Value = "  -3.55%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "386.80"
$ws.Range("E21").Value = "  -1.16%  "

# Row 22 - was Polygon, now Litecoin
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D22" "75.13"
$ws.Range("E22").Value = "  -0.11%  "

# Row 23 - was Litecoin, now Polygon
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D23" "0.560"
$ws.Range("E23").Value = "  -1.27%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.10%  "

# Row 25 - was PEPE, now WrappedeETH
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue "D25" "3.587.43"
$ws.Range("E25").Value = "  -1.14%  "

# Row 26 - was WrappedeETH, now PEPE
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D26" "0.0000115"
$ws.Range("E26").Value = "  -3.08%  "

# Row 27 - Kaspa
$ws.Range("E27").Value = "  -0.08%  "

# Row 28 - RenderToken
Set-TextValue "D28" "7.57"
$ws.Range("E28").Value = "  -3.77%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  +0.12%  "

# Row 30 - InternetComputer(DFINITY)
Set-TextValue "D30" "7.97"
$ws.Range("E30").Value = "  -4.20%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -2.40%  "

# Row 32 - USDe
$ws.Range("E32").Value = "  -0.01%  "

# Row 33 - Fetch.AI
$ws.Range("E33").Value = "  -9.92%  "

# Row 34 - EthereumClassic
Set-TextValue "D34" "23.18"
$ws.Range("E34").Value = "  -3.00%  "

# Row 35 - NEARProtocol
Set-TextValue "D35" "5.30"
$ws.Range("E35").Value = "  -1.45%  "

# Row 36 - ImmutableX
Set-TextValue "D36" "1.61"
$ws.Range("E36").Value = "  +1.89%  "

# Row 37 - EnergySwap
Set-TextValue "D37" "31.90"
$ws.Range("E37").Value = "  -0.99%  "

# Row 38 - Aptos
Set-TextValue "D38" "6.96"
$ws.Range("E38").Value = "  -2.89%  "

# Row 39 - Monero
Set-TextValue "D39" "169.97"
$ws.Range("E39").Value = "  -0.80%  "

# Row 40 - RenzoRestakedETH
Set-TextValue "D40" "3.471.88"
$ws.Range("E40").Value = "  -1.49%  "

# Row 41 - Hedera
Set-TextValue "D41" "0.0773"
$ws.Range("E41").Value = "  +0.17%  "

# Row 42 - Mantle
Set-TextValue "D42" "0.786"
$ws.Range("E42").Value = "  -3.06%  "

# Row 43 - OKB
Set-TextValue "D43" "42.47"
$ws.Range("E43").Value = "  +0.03%  "

# Row 44 - Stacks
$ws.Range("E44").Value = "  -2.62%  "

# Row 45 - Filecoin
Set-TextValue "D45" "4.33"
$ws.Range("E45").Value = "  -3.97%  "

# Row 46 - ONDO
Set-TextValue "D46" "1.17"
$ws.Range("E46").Value = "  -3.24%  "

# Row 47 - Maker
Set-TextValue "D47" "2.554.62"
$ws.Range("E47").Value = "  -2.61%  "

# Row 48 - Cosmos
Set-TextValue "D48" "6.88"
$ws.Range("E48").Value = "  +1.20%  "

# Row 49 - dogwifhat
Set-TextValue "D49" "2.25"
$ws.Range("E49").Value = "  -2.61%  "

# Row 50 - InjectiveProtocol
Set-TextValue "D50" "22.61"
$ws.Range("E50").Value = "  -4.79%  "

# Row 51 - FirstDigitalUSD
$ws.Range("E51").Value = "  +0.09%  "
